# Updates cryptos list values (price + 1h volume change) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'43.152.65"
$ws.Range('E2').Value = '  +4.52%  '
$ws.Range('D3').Value = "'2.263.91"
$ws.Range('E3').Value = '  +3.51%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'253.82"
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('D6').Value = "'0.641"
$ws.Range('E6').Value = '  +2.19%  '
$ws.Range('D7').Value = "'71.94"
$ws.Range('E7').Value = '  +5.45%  '
$ws.Range('D8').Value = "'0.677"
$ws.Range('E8').Value = '  +18.53%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').Value = "'40.20"
$ws.Range('E10').Value = '  +8.77%  '
$ws.Range('D11').Value = "'0.0976"
$ws.Range('E11').Value = '  +4.13%  '
$ws.Range('D12').Value = "'59.59"
$ws.Range('E12').Value = '  +1.23%  '
$ws.Range('D13').Value = "'7.62"
$ws.Range('E13').Value = '  +8.80%  '
$ws.Range('D14').Value = "'0.105"
$ws.Range('E14').Value = '  +1.80%  '
$ws.Range('D15').Value = "'2.604.46"
$ws.Range('E15').Value = '  +3.67%  '
$ws.Range('D16').Value = "'0.891"
$ws.Range('E16').Value = '  +2.67%  '
$ws.Range('D17').Value = "'14.88"
$ws.Range('E17').Value = '  +4.08%  '
$ws.Range('D18').Value = "'2.263.80"
$ws.Range('E18').Value = '  +2.84%  '
$ws.Range('D19').Value = "'43.022.94"
$ws.Range('E19').Value = '  +4.51%  '
$ws.Range('D20').Value = "'0.0₃0984"
$ws.Range('E20').Value = '  +2.88%  '
$ws.Range('D21').Value = "'6.31"
$ws.Range('E21').Value = '  +2.57%  '
$ws.Range('D22').Value = "'73.41"
$ws.Range('E22').Value = '  +1.90%  '
$ws.Range('D23').Value = "'237.62"
$ws.Range('E23').Value = '  +2.08%  '
$ws.Range('D24').Value = "'2.13"
$ws.Range('E24').Value = '  +5.02%  '
$ws.Range('D25').Value = "'3.90"
$ws.Range('E25').Value = '  +1.29%  '
$ws.Range('D26').Value = "'11.72"
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('D27').Value = "'0.999"
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D28').Value = "'2.47"
$ws.Range('E28').Value = '  -1.60%  '
$ws.Range('E30').Value = '  +2.21%  '
$ws.Range('D31').Value = "'168.00"
$ws.Range('E31').Value = '  -0.69%  '
$ws.Range('D32').Value = "'21.27"
$ws.Range('E32').Value = '  +3.14%  '
$ws.Range('E33').Value = '  +10.17%  '
$ws.Range('D34').Value = "'6.19"
$ws.Range('E34').Value = '  +13.60%  '
$ws.Range('D35').Value = "'0.0779"
$ws.Range('E35').Value = '  +4.77%  '
$ws.Range('E36').Value = '  +2.13%  '
$ws.Range('D37').Value = "'29.06"
$ws.Range('E37').Value = '  +11.32%  '
$ws.Range('D38').Value = "'4.75"
$ws.Range('E38').Value = '  +3.50%  '
$ws.Range('E39').Value = '  -0.51%  '
$ws.Range('D40').Value = "'0.0324"
$ws.Range('E40').Value = '  +9.25%  '
$ws.Range('D41').Value = "'2.30"
$ws.Range('E41').Value = '  +5.04%  '
$ws.Range('D42').Value = "'5.90"
$ws.Range('E42').Value = '  +4.38%  '
$ws.Range('D43').Value = "'12.46"
$ws.Range('E43').Value = '  +2.54%  '
$ws.Range('D44').Value = "'64.37"
$ws.Range('E44').Value = '  +1.72%  '
$ws.Range('E45').Value = '  +1.55%  '
$ws.Range('D46').Value = "'0.203"
$ws.Range('E46').Value = '  +2.85%  '
$ws.Range('D47').Value = "'8.96"
$ws.Range('E47').Value = '  +3.99%  '
$ws.Range('E48').Value = '  +2.49%  '
$ws.Range('E49').Value = '  -1.72%  '
$ws.Range('D50').Value = "'1.00"
$ws.Range('E50').Value = '  -0.22%  '
$ws.Range('E51').Value = '  +2.81%  '
